# Automatische test-sync: 2025-06-19 18:00:30
#
# Adds one new incoming-mail row to the "Logs" sheet, a matching tally
# row on the "Dashboard" sheet, expands the two conditional-formatting
# blocks and the bar chart's category/value series ranges to cover it.

$wb   = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Logs sheet: append row 34 ------------------------------------------
$logs.Range("A34").Value2 = "Wat zijn jullie openingstijden?"
$logs.Range("B34").Value2 = "mailmind.test@zohomail.eu"
$logs.Range("C34").Value2 = "Hallo, ik zou graag willen weten wat jullie openingstijden zijn. Dank je wel!"
$logs.Range("D34").Value2 = "Openingstijden"

# Build the multi-line reply through a formula + values-only paste so the
# embedded line breaks don't trigger Excel's row AutoFit (which would
# otherwise stamp an explicit ht="..." customHeight="1" on the new row).
$logs.Range("E34").Formula = '="Beste,"&CHAR(10)&"Bedankt voor je bericht. Onze openingstijden zijn van maandag tot en met vrijdag van 09:00 tot 18:00 uur. Op zaterdag zijn we geopend van 10:00 tot 17:00 uur. Op zondag zijn we gesloten."&CHAR(10)&"Met vriendelijke groet,"&CHAR(10)&"[Naam organisatie]"'
$logs.Range("E34").Copy()
$logs.Range("E34").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

$logs.Range("F34").Value2 = "2025-06-19 18:00:22"
$logs.Range("G34").Value2 = "Ja"

# --- expand conditional formatting to include the new row --------------
$logs.Range("D2:D33").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D34"))
$logs.Range("G2:G33").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G34"))

# --- Dashboard sheet: append row 10 -------------------------------------
$dash.Range("A10").Value2 = "Openingstijden"
$dash.Range("B10").Value2 = 1

# --- extend the chart's category/value series to row 10 ----------------
$chart = $dash.ChartObjects().Item(1).Chart
$ser = $chart.SeriesCollection().Item(1)
$ser.Formula = "=SERIES('Dashboard'!B1,'Dashboard'!`$A`$2:`$A`$10,'Dashboard'!`$B`$2:`$B`$10,1)"
